$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(0.08007198572158813,0.985791802406311,0.03856929391622543,0.9952033162117004)
    ,@(0.01875740848481655,0.9951562881469727,0.1009048670530319,0.9663468599319458)
    ,@(0.01144350972026587,0.9956299066543579,0.01072407513856888,0.998248815536499)
    ,@(0.01088115200400352,0.9955222606658936,0.006131215952336788,0.9994670152664185)
    ,@(0.008160357363522053,0.9961680769920349,0.01201025303453207,0.9987818002700806)
    ,@(0.007884222082793713,0.995608389377594,0.002876221435144544,0.9996954202651978)
    ,@(0.008362308144569397,0.9956944584846497,0.002993487054482102,0.9996193051338196)
    ,@(0.006591978948563337,0.9965125322341919,0.003327826038002968,0.9997715950012207)
    ,@(0.007527980022132397,0.996684730052948,0.005529874935746193,0.9993147850036621)
    ,@(0.007582232821732759,0.995802104473114,0.003286913502961397,0.9998477101325989)
    ,@(0.00708973640576005,0.9964479207992554,0.001985572278499603,0.9998477101325989)
    ,@(0.007011891808360815,0.9957160353660583,0.003406009636819363,0.9998477101325989)
    ,@(0.007193909492343664,0.9962111711502075,0.00235042767599225,0.9998477101325989)
    ,@(0.007274913601577282,0.9958236217498779,0.002577684121206403,0.9996954202651978)
    ,@(0.0068623386323452,0.9962111711502075,0.003303923178464174,0.9996193051338196)
    ,@(0.006823075469583273,0.9962326884269714,0.003472308395430446,0.9996954202651978)
    ,@(0.00733052846044302,0.9963833689689636,0.002757593058049679,0.9997715950012207)
    ,@(0.007271901704370975,0.9958236217498779,0.003236343618482351,0.9998477101325989)
    ,@(0.00671363016590476,0.9966632127761841,0.002509836805984378,0.9998477101325989)
    ,@(0.0066492292098701,0.9963833689689636,0.001921096933074296,0.9998477101325989)
    ,@(0.006285324692726135,0.9965340495109558,0.002470882842317224,0.9997715950012207)
    ,@(0.006875474471598864,0.9962542057037354,0.003507760586217046,0.9994670152664185)
    ,@(0.006698774173855782,0.9960173964500427,0.002723620273172855,0.9998477101325989)
    ,@(0.006915098521858454,0.9960604310035706,0.002819475950673223,0.9997715950012207)
    ,@(0.007288290187716484,0.9961035251617432,0.02501669526100159,0.9907872676849365)
    ,@(0.006761039141565561,0.996146559715271,0.002671482507139444,0.9998477101325989)
    ,@(0.006573008839040995,0.9961895942687988,0.002720019780099392,0.9998477101325989)
    ,@(0.00677180290222168,0.9967278242111206,0.002767274854704738,0.9998477101325989)
    ,@(0.006351208314299583,0.9961250424385071,0.00491533987224102,0.9996193051338196)
    ,@(0.006793205626308918,0.9964264035224915,0.008962551131844521,0.9996954202651978)
    ,@(0.006504989229142666,0.9961895942687988,0.003428101539611816,0.9996954202651978)
    ,@(0.006656272802501917,0.9960819482803345,0.003073574276641011,0.9997715950012207)
    ,@(0.006477264687418938,0.9963403344154358,0.003373962361365557,0.9997715950012207)
    ,@(0.007249389309436083,0.9962542057037354,0.004797715693712234,0.9996193051338196)
    ,@(0.006381303537636995,0.9963618516921997,0.003252320690080523,0.9997715950012207)
    ,@(0.006242933683097363,0.9965770840644836,0.003359182504937053,0.9997715950012207)
    ,@(0.006496266927570105,0.9964479207992554,0.003255567979067564,0.9997715950012207)
    ,@(0.007636864203959703,0.9961680769920349,0.007189958356320858,0.998248815536499)
    ,@(0.006333868484944105,0.9967063069343567,0.004852655809372663,0.9989340901374817)
    ,@(0.006501883268356323,0.9961680769920349,0.004438603762537241,0.9995431900024414)
    ,@(0.00640378100797534,0.9962326884269714,0.003557511139661074,0.9998477101325989)
    ,@(0.00716503756120801,0.9960389137268066,0.003508857917040586,0.9998477101325989)
    ,@(0.005996830761432648,0.9963187575340271,0.003460513660684228,0.9998477101325989)
    ,@(0.006530238781124353,0.9964694976806641,0.003481277730315924,0.9998477101325989)
    ,@(0.006652339827269316,0.9963403344154358,0.003393134335055947,0.9998477101325989)
    ,@(0.006656423676759005,0.9962326884269714,0.003866539802402258,0.9998477101325989)
    ,@(0.00738381827250123,0.9962326884269714,0.003530483925715089,0.9997715950012207)
    ,@(0.006607231218367815,0.9963403344154358,0.003549492917954922,0.9997715950012207)
    ,@(0.00632998812943697,0.9964264035224915,0.003982405178248882,0.9996193051338196)
    ,@(0.006720965728163719,0.9964048862457275,0.004111401736736298,0.9996954202651978)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt 4; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value2 = $rowVals[$j]
    }
}
